# Updates the "cryptos" price/volume table to the new scrape snapshot.
# For cells in column D whose new text is a valid numeric literal, the
# NumberFormat is forced to Text ("@") before the write so Excel keeps
# storing the cell as a string (matching the sheet's existing convention)
# instead of silently converting it to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '27.270.02'
$ws.Range("E2").Value = '  +1.21%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '1.652.22'
$ws.Range("E3").Value = '  +0.55%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.50%  '
# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.86'
$ws.Range("E5").Value = '  -0.25%  '
# Row 6: XRP
$ws.Range("E6").Value = '  +2.38%  '
# Row 7: USDC
$ws.Range("E7").Value = '  -0.52%  '
# Row 8: Cardano
$ws.Range("E8").Value = '  +2.22%  '
# Row 9: Dogecoin
$ws.Range("E9").Value = '  +0.51%  '
# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.29'
$ws.Range("E10").Value = '  +5.03%  '
# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0848'
$ws.Range("E11").Value = '  +0.02%  '
# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '1.878.74'
$ws.Range("E12").Value = '  +0.32%  '
# Row 13: WrappedEther
$ws.Range("D13").Value = '1.651.87'
$ws.Range("E13").Value = '  +0.64%  '
# Row 14: Polkadot
$ws.Range("E14").Value = '  -0.35%  '
# Row 15: Polygon
$ws.Range("E15").Value = '  +2.55%  '
# Row 16: Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.99'
$ws.Range("E16").Value = '  +3.49%  '
# Row 17: WrappedBTC
$ws.Range("D17").Value = '27.245.50'
$ws.Range("E17").Value = '  +1.20%  '
# Row 18: ShibaInu
$ws.Range("E18").Value = '  +0.98%  '
# Row 19: BitcoinCash
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '222.19'
$ws.Range("E19").Value = '  +2.21%  '
# Row 20: Dai
$ws.Range("E20").Value = '  -0.55%  '
# Row 21: Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.78'
$ws.Range("E21").Value = '  +2.86%  '
# Row 22: Uniswap
$ws.Range("E22").Value = '  +1.69%  '
# Row 23: Toncoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.48'
$ws.Range("E23").Value = '  +2.21%  '
# Row 24: Avalanche
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.26'
$ws.Range("E24").Value = '  +0.37%  '
# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.65'
$ws.Range("E25").Value = '  -0.31%  '
# Row 26: BinanceUSD
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.56%  '
# Row 27: Cosmos
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.42'
$ws.Range("E27").Value = '  +1.49%  '
# Row 28: Stellar
$ws.Range("E28").Value = '  +1.33%  '
# Row 29: EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.90'
$ws.Range("E29").Value = '  +0.38%  '
# Row 30: Hedera
$ws.Range("E30").Value = '  -0.63%  '
# Row 31: PancakeSwap
$ws.Range("E31").Value = '  -0.51%  '
# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("E32").Value = '  -0.16%  '
# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.04'
$ws.Range("E33").Value = '  +1.37%  '
# Row 34: LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.59'
$ws.Range("E34").Value = '  +1.76%  '
# Row 35: Maker
$ws.Range("D35").Value = '1.273.02'
$ws.Range("E35").Value = '  +0.37%  '
# Row 36: HuobiToken
$ws.Range("E36").Value = '  -0.17%  '
# Row 37: VeChain
$ws.Range("E37").Value = '  +3.21%  '
# Row 38: ImmutableX
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.545'
$ws.Range("E38").Value = '  +2.39%  '
# Row 39: ARBITRUM
$ws.Range("E39").Value = '  +2.93%  '
# Row 40: PaxDollar
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.53%  '
# Row 41: TrustWalletToken
$ws.Range("E41").Value = '  +0.41%  '
# Row 42: FraxShare
$ws.Range("E42").Value = '  +0.93%  '
# Row 43: MXToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.18'
$ws.Range("E43").Value = '  +4.81%  '
# Row 44: Aave
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '1.788.60'
$ws.Range("E44").Value = '  +0.33%  '
# Row 45: RocketPoolETH
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.74'
$ws.Range("E45").Value = '  +4.34%  '
# Row 46: Quant
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.46'
$ws.Range("E46").Value = '  -0.20%  '
# Row 47: RenderToken
$ws.Range("E47").Value = '  +0.37%  '
# Row 48: BabyDogeCoin
$ws.Range("E48").Value = '  +5.93%  '
# Row 49: Cronos
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0515'
$ws.Range("E49").Value = '  -0.24%  '
# Row 50: EnergySwap
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.74'
$ws.Range("E50").Value = '  +1.48%  '
# Row 51: Algorand
$ws.Range("E51").Value = '  +0.70%  '
